# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the three files that were
# "In Translation" (1afbb6e9, 245809ef, cfa930e9) get re-ordered in the
# per-language worksheets, and 1afbb6e9 moves to "Ready for handoff" status
# with a freshly stamped handoff date/time. The Overview sheet mirrors the
# same re-ordering and date stamp.

$wb = $excel.ActiveWorkbook

function Set-CellAndLink {
    param(
        $ws,
        [string]$ref,
        [string]$value,
        [string]$displayOverride = $null
    )
    $ws.Range($ref).Value = $value
    if ($ws.Range($ref).Hyperlinks.Count -ge 1) {
        $hl = $ws.Range($ref).Hyperlinks.Item(1)
        if ($displayOverride) {
            $hl.TextToDisplay = $displayOverride
        } else {
            $hl.TextToDisplay = $value
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# 1e39392e row: new handoff timestamp
$wsOverview.Range("D6").Value = "2016-21-18 14:21:36"

# Row 7 -> 245809ef (still In Translation)
Set-CellAndLink $wsOverview "A7" "245809ef-eee1-4104-a55a-7fb612c9d5d6.md"
$wsOverview.Range("B7").Value = "In Translation"
$wsOverview.Range("C7").Value = "In Translation"
$wsOverview.Range("D7").Value = "2016-20-18 14:20:59"

# Row 8 -> cfa930e9 (still In Translation)
Set-CellAndLink $wsOverview "A8" "cfa930e9-6e4a-4427-bdcf-82b107ce51c6.md"
$wsOverview.Range("B8").Value = "In Translation"
$wsOverview.Range("C8").Value = "In Translation"
$wsOverview.Range("D8").Value = "2016-15-18 14:15:46"

# Row 9 -> 1afbb6e9 (now Ready for handoff, fresh timestamp)
Set-CellAndLink $wsOverview "A9" "1afbb6e9-aa2e-4af7-9834-4be6e6cc1748.md"
$wsOverview.Range("B9").Value = "Ready for handoff"
$wsOverview.Range("C9").Value = "Ready for handoff"
$wsOverview.Range("D9").Value = "2016-21-18 14:21:36"

# Row 10 -> 1ed63e7b stays the same file but picks up the new stamp too
$wsOverview.Range("D10").Value = "2016-21-18 14:21:36"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("E6").Value = "2016-03-18 14:21:31"

Set-CellAndLink $wsZh "A7" "245809ef-eee1-4104-a55a-7fb612c9d5d6.md"
$wsZh.Range("C7").Value = "In Translation"
Set-CellAndLink $wsZh "D7" "245809ef-eee1-4104-a55a-7fb612c9d5d6.ceb3afb3e9dd3bdb7885d27c5758f9cc099895e0.zh-cn.xlf"
$wsZh.Range("E7").Value = "2016-03-18 14:20:56"

Set-CellAndLink $wsZh "A8" "cfa930e9-6e4a-4427-bdcf-82b107ce51c6.md"
$wsZh.Range("C8").Value = "In Translation"
Set-CellAndLink $wsZh "D8" "cfa930e9-6e4a-4427-bdcf-82b107ce51c6.10b79c0069adbe3264ff26409ba0fada5d329204.zh-cn.xlf"
$wsZh.Range("E8").Value = "2016-03-18 14:15:43"

Set-CellAndLink $wsZh "A9" "1afbb6e9-aa2e-4af7-9834-4be6e6cc1748.md"
$wsZh.Range("C9").Value = "Ready for handoff"
Set-CellAndLink $wsZh "D9" "1afbb6e9-aa2e-4af7-9834-4be6e6cc1748.c3c81ce89e2f83a94d242d3fd3288f3e734b3bed.zh-cn.xlf"
$wsZh.Range("E9").Value = "2016-03-18 14:21:31"

$wsZh.Range("E10").Value = "2016-03-18 14:21:31"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("E6").Value = "2016-03-18 14:21:36"

Set-CellAndLink $wsDe "A7" "245809ef-eee1-4104-a55a-7fb612c9d5d6.md"
$wsDe.Range("C7").Value = "In Translation"
Set-CellAndLink $wsDe "D7" "245809ef-eee1-4104-a55a-7fb612c9d5d6.ceb3afb3e9dd3bdb7885d27c5758f9cc099895e0.de-de.xlf"
$wsDe.Range("E7").Value = "2016-03-18 14:20:59"

Set-CellAndLink $wsDe "A8" "cfa930e9-6e4a-4427-bdcf-82b107ce51c6.md"
$wsDe.Range("C8").Value = "In Translation"
Set-CellAndLink $wsDe "D8" "cfa930e9-6e4a-4427-bdcf-82b107ce51c6.10b79c0069adbe3264ff26409ba0fada5d329204.de-de.xlf"
$wsDe.Range("E8").Value = "2016-03-18 14:15:46"

Set-CellAndLink $wsDe "A9" "1afbb6e9-aa2e-4af7-9834-4be6e6cc1748.md"
$wsDe.Range("C9").Value = "Ready for handoff"
Set-CellAndLink $wsDe "D9" "1afbb6e9-aa2e-4af7-9834-4be6e6cc1748.c3c81ce89e2f83a94d242d3fd3288f3e734b3bed.de-de.xlf"
$wsDe.Range("E9").Value = "2016-03-18 14:21:36"

$wsDe.Range("E10").Value = "2016-03-18 14:21:36"
